$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G18").Value = "Fallo"
$ws.Range("H18").Value = -1

$ws.Range("G20").Value = "Fallo"
$ws.Range("H20").Value = -1

$ws.Range("A36").Value = 14343571
$ws.Range("A37").Value = 14344407
